$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.657.82"
$ws.Range("E2").Value = "  +2.94%  "
$ws.Range("D3").Value = "'1.857.24"
$ws.Range("E3").Value = "  +1.63%  "
$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  +0.64%  "
$ws.Range("D5").Value = "'263.53"
$ws.Range("E5").Value = "  -5.24%  "
$ws.Range("E6").Value = "  +0.42%  "
$ws.Range("D7").Value = "'0.5265"
$ws.Range("E7").Value = "  +3.06%  "
$ws.Range("D8").Value = "'0.3317"
$ws.Range("E8").Value = "  -4.58%  "
$ws.Range("D9").Value = "'0.06756"
$ws.Range("E9").Value = "  -0.58%  "
$ws.Range("D10").Value = "'19.52"
$ws.Range("E10").Value = "  -1.86%  "
$ws.Range("D11").Value = "'0.7772"
$ws.Range("E11").Value = "  -4.02%  "
$ws.Range("D12").Value = "'0.07738"
$ws.Range("E12").Value = "  -0.96%  "
$ws.Range("D13").Value = "'1.830.62"
$ws.Range("E13").Value = "  +0.20%  "
$ws.Range("D14").Value = "'89.42"
$ws.Range("E14").Value = "  +1.58%  "
$ws.Range("D15").Value = "'5.078"
$ws.Range("E15").Value = "  +0.05%  "
$ws.Range("D16").Value = "'1.005"
$ws.Range("E16").Value = "  +0.83%  "
$ws.Range("D17").Value = "'14.28"
$ws.Range("E17").Value = "  +0.89%  "
$ws.Range("D18").Value = "'1.002"
$ws.Range("E18").Value = "  +0.35%  "
$ws.Range("D19").Value = "'0.000007935"
$ws.Range("E19").Value = "  -1.60%  "
$ws.Range("D20").Value = "'26.698.48"
$ws.Range("E20").Value = "  +2.95%  "
$ws.Range("D21").Value = "'2.062.92"
$ws.Range("E21").Value = "  +0.96%  "
$ws.Range("D22").Value = "'4.629"
$ws.Range("E22").Value = "  -2.91%  "
$ws.Range("D23").Value = "'9.796"
$ws.Range("E23").Value = "  -2.23%  "
$ws.Range("D24").Value = "'6.004"
$ws.Range("E24").Value = "  -2.94%  "
$ws.Range("D25").Value = "'2.382"
$ws.Range("E25").Value = "  +1.05%  "
$ws.Range("D26").Value = "'145.64"
$ws.Range("E26").Value = "  +1.85%  "
$ws.Range("D27").Value = "'1.664"
$ws.Range("E27").Value = "  -0.04%  "
$ws.Range("D28").Value = "'17.11"
$ws.Range("E28").Value = "  -0.50%  "
$ws.Range("D29").Value = "'112.25"
$ws.Range("E29").Value = "  +2.63%  "
$ws.Range("D30").Value = "'4.254"
$ws.Range("E30").Value = "  -1.84%  "
$ws.Range("D31").Value = "'4.237"
$ws.Range("E31").Value = "  -1.34%  "
$ws.Range("D32").Value = "'0.08799"
$ws.Range("E32").Value = "  +0.57%  "
$ws.Range("D33").Value = "'0.04914"
$ws.Range("E33").Value = "  +1.04%  "
$ws.Range("D34").Value = "'1.147"
$ws.Range("E34").Value = "  -1.79%  "
$ws.Range("D35").Value = "'2.883"
$ws.Range("E35").Value = "  +0.85%  "
$ws.Range("D36").Value = "'0.7154"
$ws.Range("E36").Value = "  -1.98%  "
$ws.Range("D37").Value = "'3.152"
$ws.Range("E37").Value = "  -0.93%  "
$ws.Range("D38").Value = "'0.01825"
$ws.Range("E38").Value = "  -1.46%  "
$ws.Range("D39").Value = "'2.258"
$ws.Range("E39").Value = "  -5.91%  "
$ws.Range("D40").Value = "'0.4996"
$ws.Range("E40").Value = "  -2.48%  "
$ws.Range("D41").Value = "'115.56"
$ws.Range("E41").Value = "  -1.18%  "
$ws.Range("D42").Value = "'0.9103"
$ws.Range("E42").Value = "  -3.99%  "
$ws.Range("D43").Value = "'6.088"
$ws.Range("E43").Value = "  -1.93%  "
$ws.Range("D44").Value = "'7.841"
$ws.Range("E44").Value = "  -2.11%  "
$ws.Range("D45").Value = "'1.002"
$ws.Range("E45").Value = "  +0.45%  "
$ws.Range("D46").Value = "'0.4320"
$ws.Range("E46").Value = "  -3.79%  "
$ws.Range("D47").Value = "'0.1306"
$ws.Range("E47").Value = "  -4.19%  "
$ws.Range("D48").Value = "'9.208"
$ws.Range("E48").Value = "  -0.69%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "'0.05930"
$ws.Range("E49").Value = "  +0.26%  "
$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").Value = "'35.69"
$ws.Range("E50").Value = "  -1.53%  "
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "'60.26"
$ws.Range("E51").Value = "  -0.53%  "
